# Updates the "Problem Statement" paragraphs so the wording switches from
# second-person ("you/your") to first-person-plural ("we/our") phrasing, per
# the commit diff. The rest of the document's visible text is unchanged.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "When you have a great business idea, funding is nearly always the sticking point. It" + [char]0x2019 + "s a great idea, after all, but how can you raise the money to get it started?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "When we have a great business idea, funding is nearly always the sticking point. It" + [char]0x2019 + "s a great idea, after all, but how can we raise the money to get it started?",
    2
)

$d.Content.Find.Execute(
    "If you have a tech-based idea, you may have an easier time attracting attention from venture capitalists or angel investors, but as more companies work that angle, finding an investor is harder than ever. So how can you get your business off the ground?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "If we have a tech-based idea, we may have an easier time attracting attention from venture capitalists or angel investors, but as more companies work that angle, finding an investor is harder than ever. So how can we get our business off the ground?",
    2
)

$d.Content.Find.Execute(
    "The downside? They are often focused on tech-heavy businesses, so you might struggle to find one that works for your company.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The downside? They are often focused on tech-heavy businesses, so we might struggle to find one that works for our company.",
    2
)
